$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129 (this shifts the existing rows 129-154 down to 130-155,
# and keeps formatting/styles consistent with the surrounding rows).
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly price entry.
$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44476
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = 100112043
$ws.Range("G129").Value = "Pepino ensalada"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 160
$ws.Range("K129").Value = 22000
$ws.Range("L129").Value = 22000
$ws.Range("M129").Value = 22000
$ws.Range("N129").Value = "$/caja 60 unidades"
$ws.Range("O129").Value = "Región de Arica y Parinacota"
$ws.Range("P129").Value = 367
$ws.Range("Q129").Value = 60
$ws.Range("R129").Value = "Hortaliza"
